$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.109881997108459
$ws.Range("B1").Value = 2.243576765060425
$ws.Range("C1").Value = 10.11110782623291
$ws.Range("D1").Value = 1.439636826515198
$ws.Range("E1").Value = 1.287283897399902
